# Update "想去人数" (people interested) counts in column F across sheets.
# Mirrors a re-scrape of the source data published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3150
$ws.Range("F4").Value = 1964
$ws.Range("F5").Value = 258
$ws.Range("F6").Value = 84
$ws.Range("F7").Value = 2627
$ws.Range("F8").Value = 599
$ws.Range("F9").Value = 297
$ws.Range("F10").Value = 29
$ws.Range("F13").Value = 142
$ws.Range("F14").Value = 9909
$ws.Range("F15").Value = 66
$ws.Range("F17").Value = 2
$ws.Range("F19").Value = 7819
$ws.Range("F20").Value = 12426
$ws.Range("F24").Value = 391
$ws.Range("F25").Value = 588
$ws.Range("F26").Value = 2786
$ws.Range("F27").Value = 253
$ws.Range("F28").Value = 229
$ws.Range("F29").Value = 7838
$ws.Range("F30").Value = 1350
$ws.Range("F32").Value = 69
$ws.Range("F33").Value = 73
$ws.Range("F34").Value = 4599
$ws.Range("F35").Value = 1290
$ws.Range("F36").Value = 62
$ws.Range("F37").Value = 371
$ws.Range("F39").Value = 612

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 70
$ws.Range("F15").Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 3150
$ws.Range("F6").Value = 1964
$ws.Range("F8").Value = 258
$ws.Range("F9").Value = 2627
$ws.Range("F11").Value = 599
$ws.Range("F12").Value = 297
$ws.Range("F13").Value = 29
$ws.Range("F16").Value = 142
$ws.Range("F17").Value = 9909
$ws.Range("F18").Value = 66
$ws.Range("F21").Value = 7819
$ws.Range("F22").Value = 12427
$ws.Range("F27").Value = 588
$ws.Range("F29").Value = 2786
$ws.Range("F32").Value = 253
$ws.Range("F33").Value = 229
$ws.Range("F35").Value = 69
$ws.Range("F36").Value = 73
$ws.Range("F37").Value = 4601
$ws.Range("F38").Value = 70
$ws.Range("F40").Value = 11
$ws.Range("F45").Value = 612
